# Restore Rule R30's "Integer min" threshold (cell C10 on the "Rules" sheet)
# from 18 to 1, per revision #c7d9c6270511395650ebf58b802cbfc4ee6f18d9.TEST

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
